$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the N3J6 match-day data (minutes played in DS, T/R/NR/HG status in DT,
# goals in DU, assists in DV) for each player row.
$ws.Range("DS2").Value = 90
$ws.Range("DT2").Value = "T"
$ws.Range("DT3").Value = "HG"
$ws.Range("DT4").Value = "NR"
$ws.Range("DS5").Value = 90
$ws.Range("DT5").Value = "T"
$ws.Range("DT6").Value = "HG"
$ws.Range("DS7").Value = 90
$ws.Range("DT7").Value = "T"
$ws.Range("DT8").Value = "HG"
$ws.Range("DS9").Value = 90
$ws.Range("DT9").Value = "T"
$ws.Range("DT10").Value = "HG"
$ws.Range("DT11").Value = "NR"
$ws.Range("DT12").Value = "HG"
$ws.Range("DT13").Value = "HG"
$ws.Range("DS14").Value = 90
$ws.Range("DT14").Value = "T"
$ws.Range("DU14").Value = 1
$ws.Range("DS15").Value = 60
$ws.Range("DT15").Value = "T"
$ws.Range("DS16").Value = 88
$ws.Range("DT16").Value = "T"
$ws.Range("DV16").Value = 1
$ws.Range("DT17").Value = "HG"
$ws.Range("DS18").Value = 30
$ws.Range("DT18").Value = "R"
$ws.Range("DS19").Value = 2
$ws.Range("DT19").Value = "R"
$ws.Range("DS20").Value = 80
$ws.Range("DT20").Value = "T"
$ws.Range("DT21").Value = "HG"
$ws.Range("DS22").Value = 90
$ws.Range("DT22").Value = "T"
$ws.Range("DT23").Value = "HG"
$ws.Range("DT24").Value = "HG"
$ws.Range("DS25").Value = 10
$ws.Range("DT25").Value = "R"
$ws.Range("DT26").Value = "HG"
$ws.Range("DS27").Value = 90
$ws.Range("DT27").Value = "T"
$ws.Range("DS28").Value = 90
$ws.Range("DT28").Value = "T"

# Update the view to reflect the scroll/selection position left after the edit.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("DR1").Select()
$win.FreezePanes = $true
$ws.Range("DX23").Select()
